$d = $word.ActiveDocument

# --- Paragraph 1: "To: " -> "To: APPLIED MEDICAL RESOURCES" (plain, unformatted run) ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$r1.Delete()
$ins1 = $d.Range($p1.Range.Start, $p1.Range.Start)
$ins1.InsertBefore("To: APPLIED MEDICAL RESOURCES")

# --- Paragraph 4: "Date: " -> "Date:  2022-12-30" (plain, unformatted run) ---
$p4 = $d.Paragraphs.Item(4)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$r4.Delete()
$ins4 = $d.Range($p4.Range.Start, $p4.Range.Start)
$ins4.InsertBefore("Date:  2022-12-30")
